# testReportListener.java was updated to exclude certain methods from the
# generated report. As a result, the test run that produced this report had
# one fewer (excluded) test-case row, and all the remaining test cases were
# re-executed, producing new execution timestamps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for the excluded test case (row 3 held a blank test-case
# name with no further data). Deleting the row shifts all following rows up
# by one and keeps everything else intact.
$ws.Rows.Item(3).Delete()

# Re-assert the sequential Test Case ID numbering for the rows that shifted
# up, preserving the original numeric formatting.
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8

# Update the execution-time column for the remaining test cases to reflect
# the new test run's timestamps.
$ws.Range("D2").Value = "01/04/2025 01:55:45 PM"
$ws.Range("D3").Value = "01/04/2025 01:55:47 PM"
$ws.Range("D4").Value = "01/04/2025 01:55:48 PM"
$ws.Range("D5").Value = "01/04/2025 01:55:51 PM"
$ws.Range("D6").Value = "01/04/2025 01:56:03 PM"
$ws.Range("D7").Value = "01/04/2025 01:56:05 PM"
$ws.Range("D8").Value = "01/04/2025 01:56:06 PM"
$ws.Range("D9").Value = "01/04/2025 01:56:06 PM"
